$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new sheet "2022-Q3" right before "2022-Q1" (i.e. right after
#    the "总计" summary sheet), and fill it with the quarterly fund data.
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q1Sheet      = $wb.Worksheets.Item(2)

$newSheet = $wb.Worksheets.Add($q1Sheet)
$newSheet.Name = "2022-Q3"

# Header row (B1:H1) — reuse the bold/centered/bordered style already used
# for header cells on the "总计" sheet.
$summarySheet.Cells.Item(1, 2).Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Index column (A2:A3) — reuse the same index-column style.
$summarySheet.Cells.Item(2, 1).Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(3, 1).Value = 1

# Row 2: 华夏博锐一年持有混合（MOM）A
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Cells.Item(2, 2).Value = "011361"
$newSheet.Cells.Item(2, 3).Value = "华夏博锐一年持有混合（MOM）A"
$newSheet.Cells.Item(2, 4).Value = "0.11"
$newSheet.Cells.Item(2, 5).Value = "33.21"
$newSheet.Cells.Item(2, 6).Value = "1.23"
$newSheet.Cells.Item(2, 7).Value = "0.0014"
$newSheet.Range("B2:G2").Style = "Normal"
$newSheet.Cells.Item(2, 8).Value = 10

# Row 3: 华夏博锐一年持有混合（MOM）C
$newSheet.Range("B3:F3").NumberFormat = "@"
$newSheet.Cells.Item(3, 2).Value = "011362"
$newSheet.Cells.Item(3, 3).Value = "华夏博锐一年持有混合（MOM）C"
$newSheet.Cells.Item(3, 4).Value = "0.00"
$newSheet.Cells.Item(3, 5).Value = "33.21"
$newSheet.Cells.Item(3, 6).Value = "1.23"
$newSheet.Range("B3:F3").Style = "Normal"
$newSheet.Cells.Item(3, 7).Value = 0
$newSheet.Cells.Item(3, 8).Value = 10

# ---------------------------------------------------------------------------
# 2. Insert a new row into the "总计" sheet for the 2022-Q3 totals, pushing
#    the existing rows down and bumping their running index in column A.
# ---------------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()
$summarySheet.Range("B2:D2").ClearFormats()

$summarySheet.Cells.Item(3, 1).Copy()
$summarySheet.Cells.Item(2, 1).PasteSpecial(-4122)

$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 2
$summarySheet.Cells.Item(2, 4).Value = 0

$summarySheet.Cells.Item(3, 1).Value = 1
$summarySheet.Cells.Item(4, 1).Value = 2
$summarySheet.Cells.Item(5, 1).Value = 3
